$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 0.6659359931945801
$ws.Cells.Item(2, 5).Value = 256.0928999514053
$ws.Cells.Item(2, 6).Value = 0.007032443523326081
$ws.Cells.Item(2, 7).Value = 0.006783430450829687
$ws.Cells.Item(2, 8).Value = 0.006474159709730734
$ws.Cells.Item(2, 9).Value = 0.005950171742866093
$ws.Cells.Item(2, 10).Value = 0.005950171742866093
$ws.Cells.Item(2, 11).Value = 0.005922025132729625
$ws.Cells.Item(2, 12).Value = 0.00564876393746221
$ws.Cells.Item(2, 13).Value = 0.00564876393746221
$ws.Cells.Item(2, 14).Value = 0.005442130296274736
$ws.Cells.Item(2, 15).Value = 0.005256475974481298
$ws.Cells.Item(2, 16).Value = 0.005256475974481298
$ws.Cells.Item(2, 17).Value = 0.005256475974481298
$ws.Cells.Item(2, 18).Value = 0.005229287383972085
$ws.Cells.Item(2, 19).Value = 0.005184028522673567
$ws.Cells.Item(2, 20).Value = 0.00516040457227762
$ws.Cells.Item(2, 21).Value = 0.005123348650963838
$ws.Cells.Item(2, 22).Value = 0.005036594366590248
$ws.Cells.Item(2, 23).Value = 0.005036594366590248
$ws.Cells.Item(2, 24).Value = 0.005012193982238953
$ws.Cells.Item(2, 25).Value = 0.004992064326538115

$ws.Cells.Item(3, 3).Value = 0.8569977283477783
$ws.Cells.Item(3, 5).Value = 255.5170934523812
$ws.Cells.Item(3, 6).Value = 0.007081895806242496
$ws.Cells.Item(3, 7).Value = 0.006327455916786053
$ws.Cells.Item(3, 8).Value = 0.006097451413930315
$ws.Cells.Item(3, 9).Value = 0.006097451413930315
$ws.Cells.Item(3, 10).Value = 0.006008958446964809
$ws.Cells.Item(3, 11).Value = 0.005932570280577505
$ws.Cells.Item(3, 12).Value = 0.005634578308541413
$ws.Cells.Item(3, 13).Value = 0.005624364626148472
$ws.Cells.Item(3, 14).Value = 0.005412309996340338
$ws.Cells.Item(3, 15).Value = 0.005412309996340338
$ws.Cells.Item(3, 16).Value = 0.005328884184363396
$ws.Cells.Item(3, 17).Value = 0.005186904657414585
$ws.Cells.Item(3, 18).Value = 0.005135670864199436
$ws.Cells.Item(3, 19).Value = 0.005061836404251042
$ws.Cells.Item(3, 20).Value = 0.005061836404251042
$ws.Cells.Item(3, 21).Value = 0.005061836404251042
$ws.Cells.Item(3, 22).Value = 0.004998571819327564
$ws.Cells.Item(3, 23).Value = 0.004998571819327564
$ws.Cells.Item(3, 24).Value = 0.004992888889777368
$ws.Cells.Item(3, 25).Value = 0.004980840028311524

$ws.Cells.Item(4, 3).Value = 0.9159939289093018
$ws.Cells.Item(4, 5).Value = 256.4854172430896
$ws.Cells.Item(4, 6).Value = 0.0073219678944489
$ws.Cells.Item(4, 7).Value = 0.006556643341068973
$ws.Cells.Item(4, 8).Value = 0.00609967631307222
$ws.Cells.Item(4, 9).Value = 0.00604630999794328
$ws.Cells.Item(4, 10).Value = 0.005852989121622069
$ws.Cells.Item(4, 11).Value = 0.005816090946170862
$ws.Cells.Item(4, 12).Value = 0.005733659724177356
$ws.Cells.Item(4, 13).Value = 0.005647068763278086
$ws.Cells.Item(4, 14).Value = 0.005337413191583359
$ws.Cells.Item(4, 15).Value = 0.005289565053378724
$ws.Cells.Item(4, 16).Value = 0.005289565053378724
$ws.Cells.Item(4, 17).Value = 0.005273674216377078
$ws.Cells.Item(4, 18).Value = 0.005219019832679973
$ws.Cells.Item(4, 19).Value = 0.005219019832679973
$ws.Cells.Item(4, 20).Value = 0.005130417540501484
$ws.Cells.Item(4, 21).Value = 0.005071791721611961
$ws.Cells.Item(4, 22).Value = 0.005071791721611961
$ws.Cells.Item(4, 23).Value = 0.005034409150629245
$ws.Cells.Item(4, 24).Value = 0.005024283830554944
$ws.Cells.Item(4, 25).Value = 0.00499971573573274

$ws.Cells.Item(5, 3).Value = 0.7849998474121094
$ws.Cells.Item(5, 5).Value = 267.5803483547861
$ws.Cells.Item(5, 6).Value = 0.007486353350393142
$ws.Cells.Item(5, 7).Value = 0.006351103317569703
$ws.Cells.Item(5, 8).Value = 0.006335741669135153
$ws.Cells.Item(5, 9).Value = 0.006093752126787336
$ws.Cells.Item(5, 10).Value = 0.006061433561310789
$ws.Cells.Item(5, 11).Value = 0.006020543586300189
$ws.Cells.Item(5, 12).Value = 0.006020543586300189
$ws.Cells.Item(5, 13).Value = 0.005897494833714598
$ws.Cells.Item(5, 14).Value = 0.005866717293662665
$ws.Cells.Item(5, 15).Value = 0.005720222348196935
$ws.Cells.Item(5, 16).Value = 0.005720222348196935
$ws.Cells.Item(5, 17).Value = 0.005599674086144204
$ws.Cells.Item(5, 18).Value = 0.005453884060557241
$ws.Cells.Item(5, 19).Value = 0.005389324888215448
$ws.Cells.Item(5, 20).Value = 0.005385956433841336
$ws.Cells.Item(5, 21).Value = 0.005347666516843152
$ws.Cells.Item(5, 22).Value = 0.005347666516843152
$ws.Cells.Item(5, 23).Value = 0.005259560575071121
$ws.Cells.Item(5, 24).Value = 0.005242049179430244
$ws.Cells.Item(5, 25).Value = 0.005215991195999728

$ws.Cells.Item(6, 3).Value = 0.7859978675842285
$ws.Cells.Item(6, 5).Value = 260.2497147454833
$ws.Cells.Item(6, 6).Value = 0.007403844270407595
$ws.Cells.Item(6, 7).Value = 0.006427717238517362
$ws.Cells.Item(6, 8).Value = 0.006427717238517362
$ws.Cells.Item(6, 9).Value = 0.006224734913735801
$ws.Cells.Item(6, 10).Value = 0.005945923442596344
$ws.Cells.Item(6, 11).Value = 0.005945923442596344
$ws.Cells.Item(6, 12).Value = 0.005749836103082327
$ws.Cells.Item(6, 13).Value = 0.005696897106285865
$ws.Cells.Item(6, 14).Value = 0.005643435606415573
$ws.Cells.Item(6, 15).Value = 0.005472900252861876
$ws.Cells.Item(6, 16).Value = 0.005439817593585005
$ws.Cells.Item(6, 17).Value = 0.005370126879223351
$ws.Cells.Item(6, 18).Value = 0.005364682069529199
$ws.Cells.Item(6, 19).Value = 0.005333050028311915
$ws.Cells.Item(6, 20).Value = 0.005148625743973507
$ws.Cells.Item(6, 21).Value = 0.005148625743973507
$ws.Cells.Item(6, 22).Value = 0.005135552355175774
$ws.Cells.Item(6, 23).Value = 0.00512772606758932
$ws.Cells.Item(6, 24).Value = 0.005087794799642102
$ws.Cells.Item(6, 25).Value = 0.005073093854687782

$ws.Cells.Item(7, 3).Value = 0.793992280960083
$ws.Cells.Item(7, 5).Value = 265.8791439196666
$ws.Cells.Item(7, 6).Value = 0.007421967774258326
$ws.Cells.Item(7, 7).Value = 0.006671704039495718
$ws.Cells.Item(7, 8).Value = 0.006429937046541276
$ws.Cells.Item(7, 9).Value = 0.006401503361879551
$ws.Cells.Item(7, 10).Value = 0.00577805443705925
$ws.Cells.Item(7, 11).Value = 0.00577805443705925
$ws.Cells.Item(7, 12).Value = 0.00577805443705925
$ws.Cells.Item(7, 13).Value = 0.00577805443705925
$ws.Cells.Item(7, 14).Value = 0.005755870272946709
$ws.Cells.Item(7, 15).Value = 0.005713637858736713
$ws.Cells.Item(7, 16).Value = 0.005635567739202413
$ws.Cells.Item(7, 17).Value = 0.005495666250457179
$ws.Cells.Item(7, 18).Value = 0.00545477309040326
$ws.Cells.Item(7, 19).Value = 0.005293673378196934
$ws.Cells.Item(7, 20).Value = 0.005275641820053687
$ws.Cells.Item(7, 21).Value = 0.005275641820053687
$ws.Cells.Item(7, 22).Value = 0.005222732130278123
$ws.Cells.Item(7, 23).Value = 0.005199816786593904
$ws.Cells.Item(7, 24).Value = 0.005182829316172838
$ws.Cells.Item(7, 25).Value = 0.005182829316172838

$ws.Cells.Item(8, 3).Value = 0.8070001602172852
$ws.Cells.Item(8, 5).Value = 263.7722736083288
$ws.Cells.Item(8, 6).Value = 0.007490561769204342
$ws.Cells.Item(8, 7).Value = 0.006484432836190194
$ws.Cells.Item(8, 8).Value = 0.00613197924777284
$ws.Cells.Item(8, 9).Value = 0.006062214963209294
$ws.Cells.Item(8, 10).Value = 0.006062214963209294
$ws.Cells.Item(8, 11).Value = 0.005912678506296744
$ws.Cells.Item(8, 12).Value = 0.005423562433509244
$ws.Cells.Item(8, 13).Value = 0.005374424008866455
$ws.Cells.Item(8, 14).Value = 0.005374424008866455
$ws.Cells.Item(8, 15).Value = 0.005373220289775895
$ws.Cells.Item(8, 16).Value = 0.005263467639606513
$ws.Cells.Item(8, 17).Value = 0.005263467639606513
$ws.Cells.Item(8, 18).Value = 0.005263467639606513
$ws.Cells.Item(8, 19).Value = 0.005263467639606513
$ws.Cells.Item(8, 20).Value = 0.005236132538354075
$ws.Cells.Item(8, 21).Value = 0.005194149552503465
$ws.Cells.Item(8, 22).Value = 0.005163375876228869
$ws.Cells.Item(8, 23).Value = 0.005163375876228869
$ws.Cells.Item(8, 24).Value = 0.005154773248188156
$ws.Cells.Item(8, 25).Value = 0.005141759719460598

$ws.Cells.Item(9, 3).Value = 0.7019977569580078
$ws.Cells.Item(9, 5).Value = 253.8378651168732
$ws.Cells.Item(9, 6).Value = 0.00735617755529799
$ws.Cells.Item(9, 7).Value = 0.006360479978485533
$ws.Cells.Item(9, 8).Value = 0.00615620542904446
$ws.Cells.Item(9, 9).Value = 0.005948897210727538
$ws.Cells.Item(9, 10).Value = 0.005776401591559044
$ws.Cells.Item(9, 11).Value = 0.005489265868975556
$ws.Cells.Item(9, 12).Value = 0.005489265868975556
$ws.Cells.Item(9, 13).Value = 0.005397742013480422
$ws.Cells.Item(9, 14).Value = 0.005350943602394113
$ws.Cells.Item(9, 15).Value = 0.005350943602394113
$ws.Cells.Item(9, 16).Value = 0.005250149608096518
$ws.Cells.Item(9, 17).Value = 0.005175130531738884
$ws.Cells.Item(9, 18).Value = 0.005175130531738884
$ws.Cells.Item(9, 19).Value = 0.005100857896939994
$ws.Cells.Item(9, 20).Value = 0.005074982563896766
$ws.Cells.Item(9, 21).Value = 0.005044837611496995
$ws.Cells.Item(9, 22).Value = 0.005036218784536709
$ws.Cells.Item(9, 23).Value = 0.004992303382959573
$ws.Cells.Item(9, 24).Value = 0.004969398306123173
$ws.Cells.Item(9, 25).Value = 0.004948106532492654

$ws.Cells.Item(10, 3).Value = 0.7949802875518799
$ws.Cells.Item(10, 5).Value = 257.3475159033242
$ws.Cells.Item(10, 6).Value = 0.007025550910566343
$ws.Cells.Item(10, 7).Value = 0.006569124012019112
$ws.Cells.Item(10, 8).Value = 0.006378595821856168
$ws.Cells.Item(10, 9).Value = 0.006058313738248416
$ws.Cells.Item(10, 10).Value = 0.005598021288323221
$ws.Cells.Item(10, 11).Value = 0.005598021288323221
$ws.Cells.Item(10, 12).Value = 0.005495130997239496
$ws.Cells.Item(10, 13).Value = 0.005385475261446648
$ws.Cells.Item(10, 14).Value = 0.005385475261446648
$ws.Cells.Item(10, 15).Value = 0.005305062349539582
$ws.Cells.Item(10, 16).Value = 0.005305062349539582
$ws.Cells.Item(10, 17).Value = 0.005257357412179951
$ws.Cells.Item(10, 18).Value = 0.00525715317440522
$ws.Cells.Item(10, 19).Value = 0.005224268284806636
$ws.Cells.Item(10, 20).Value = 0.00515173135752197
$ws.Cells.Item(10, 21).Value = 0.00515173135752197
$ws.Cells.Item(10, 22).Value = 0.00511612895881275
$ws.Cells.Item(10, 23).Value = 0.005051176026043828
$ws.Cells.Item(10, 24).Value = 0.005040913059806078
$ws.Cells.Item(10, 25).Value = 0.005016520777842575

$ws.Cells.Item(11, 3).Value = 0.7430274486541748
$ws.Cells.Item(11, 5).Value = 263.0978521399502
$ws.Cells.Item(11, 6).Value = 0.007213681725092166
$ws.Cells.Item(11, 7).Value = 0.006395097112607309
$ws.Cells.Item(11, 8).Value = 0.006381199415738899
$ws.Cells.Item(11, 9).Value = 0.006178138722840889
$ws.Cells.Item(11, 10).Value = 0.005936587801741579
$ws.Cells.Item(11, 11).Value = 0.005936587801741579
$ws.Cells.Item(11, 12).Value = 0.005936587801741579
$ws.Cells.Item(11, 13).Value = 0.005757621510997638
$ws.Cells.Item(11, 14).Value = 0.005724808862069464
$ws.Cells.Item(11, 15).Value = 0.005590992112839693
$ws.Cells.Item(11, 16).Value = 0.005590992112839693
$ws.Cells.Item(11, 17).Value = 0.005473955577130587
$ws.Cells.Item(11, 18).Value = 0.005444397725931758
$ws.Cells.Item(11, 19).Value = 0.00532282333466068
$ws.Cells.Item(11, 20).Value = 0.00525424208816154
$ws.Cells.Item(11, 21).Value = 0.00522117490343951
$ws.Cells.Item(11, 22).Value = 0.005176866429845576
$ws.Cells.Item(11, 23).Value = 0.005176866429845576
$ws.Cells.Item(11, 24).Value = 0.005134067571658415
$ws.Cells.Item(11, 25).Value = 0.005128613102143277
